$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("IMX-USD", "TAO-USD", "MNT-USD", "GRT-USD")

$startRow = 322
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
